# The daily driver report rows (2-16) were reshuffled into a new order.
# Row 8 (James Davis / TRK-1005) keeps its original position and content.
# Every other row's full record (name, asset, site, times, status, reason)
# was moved as a whole to a new row. We rewrite each data row in place,
# cell by cell, with the values from the reordered record so the final
# sheet matches the target layout exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Mark Thompson'
$ws.Range("B2").Value = 'TRK-1015'
$ws.Range("C2").Value = 'Downtown Construction'
$ws.Range("D2").Value = '06:00 AM'
$ws.Range("E2").Value = '03:00 PM'
$ws.Range("F2").Value = 'N/A'
$ws.Range("G2").Value = 'N/A'
$ws.Range("H2").Value = 'On Time'
$ws.Range("I2").Value = ''

$ws.Range("A3").Value = 'Christopher Thomas'
$ws.Range("B3").Value = 'TRK-1011'
$ws.Range("C3").Value = 'Downtown Construction'
$ws.Range("D3").Value = '06:30 AM'
$ws.Range("E3").Value = '03:30 PM'
$ws.Range("F3").Value = '06:29 AM'
$ws.Range("G3").Value = '03:33 PM'
$ws.Range("H3").Value = 'On Time'
$ws.Range("I3").Value = ''

$ws.Range("A4").Value = 'John Smith'
$ws.Range("B4").Value = 'TRK-1001'
$ws.Range("C4").Value = 'North Dallas Site'
$ws.Range("D4").Value = '06:30 AM'
$ws.Range("E4").Value = '03:30 PM'
$ws.Range("F4").Value = '06:25 AM'
$ws.Range("G4").Value = '03:35 PM'
$ws.Range("H4").Value = 'On Time'
$ws.Range("I4").Value = ''

$ws.Range("A5").Value = 'Richard Wilson'
$ws.Range("B5").Value = 'TRK-1007'
$ws.Range("C5").Value = 'Downtown Construction'
$ws.Range("D5").Value = '06:15 AM'
$ws.Range("E5").Value = '03:15 PM'
$ws.Range("F5").Value = '06:14 AM'
$ws.Range("G5").Value = '03:18 PM'
$ws.Range("H5").Value = 'On Time'
$ws.Range("I5").Value = ''

$ws.Range("A6").Value = 'William Brown'
$ws.Range("B6").Value = 'TRK-1004'
$ws.Range("C6").Value = 'Richardson Development'
$ws.Range("D6").Value = '06:30 AM'
$ws.Range("E6").Value = '03:30 PM'
$ws.Range("F6").Value = '06:32 AM'
$ws.Range("G6").Value = '02:54 PM'
$ws.Range("H6").Value = 'Early End'
$ws.Range("I6").Value = '36 minutes early'

$ws.Range("A7").Value = 'Thomas Taylor'
$ws.Range("B7").Value = 'TRK-1009'
$ws.Range("C7").Value = 'North Dallas Site'
$ws.Range("D7").Value = '07:00 AM'
$ws.Range("E7").Value = '04:00 PM'
$ws.Range("F7").Value = '06:43 AM'
$ws.Range("G7").Value = '03:52 PM'
$ws.Range("H7").Value = 'On Time'
$ws.Range("I7").Value = ''

$ws.Range("A9").Value = 'Joseph Moore'
$ws.Range("B9").Value = 'TRK-1008'
$ws.Range("C9").Value = 'Richardson Development'
$ws.Range("D9").Value = '06:30 AM'
$ws.Range("E9").Value = '03:30 PM'
$ws.Range("F9").Value = '06:55 AM'
$ws.Range("G9").Value = '03:25 PM'
$ws.Range("H9").Value = 'Late'
$ws.Range("I9").Value = '25 minutes late'

$ws.Range("A10").Value = 'Matthew Harris'
$ws.Range("B10").Value = 'TRK-1013'
$ws.Range("C10").Value = 'North Dallas Site'
$ws.Range("D10").Value = '06:45 AM'
$ws.Range("E10").Value = '03:45 PM'
$ws.Range("F10").Value = '06:52 AM'
$ws.Range("G10").Value = '03:37 PM'
$ws.Range("H10").Value = 'On Time'
$ws.Range("I10").Value = ''

$ws.Range("A11").Value = 'David Miller'
$ws.Range("B11").Value = 'TRK-1006'
$ws.Range("C11").Value = 'West Plano Project'
$ws.Range("D11").Value = '06:45 AM'
$ws.Range("E11").Value = '03:45 PM'
$ws.Range("F11").Value = '06:47 AM'
$ws.Range("G11").Value = '03:42 PM'
$ws.Range("H11").Value = 'On Time'
$ws.Range("I11").Value = ''

$ws.Range("A12").Value = 'Michael Johnson'
$ws.Range("B12").Value = 'TRK-1002'
$ws.Range("C12").Value = 'West Plano Project'
$ws.Range("D12").Value = '06:00 AM'
$ws.Range("E12").Value = '03:00 PM'
$ws.Range("F12").Value = '06:18 AM'
$ws.Range("G12").Value = '03:07 PM'
$ws.Range("H12").Value = 'Late'
$ws.Range("I12").Value = '18 minutes late'

$ws.Range("A13").Value = 'Daniel Jackson'
$ws.Range("B13").Value = 'TRK-1012'
$ws.Range("C13").Value = 'Richardson Development'
$ws.Range("D13").Value = '06:15 AM'
$ws.Range("E13").Value = '03:15 PM'
$ws.Range("F13").Value = '06:16 AM'
$ws.Range("G13").Value = '03:12 PM'
$ws.Range("H13").Value = 'On Time'
$ws.Range("I13").Value = ''

$ws.Range("A14").Value = 'Robert Williams'
$ws.Range("B14").Value = 'TRK-1003'
$ws.Range("C14").Value = 'Downtown Construction'
$ws.Range("D14").Value = '07:00 AM'
$ws.Range("E14").Value = '04:00 PM'
$ws.Range("F14").Value = '07:28 AM'
$ws.Range("G14").Value = '03:45 PM'
$ws.Range("H14").Value = 'Late'
$ws.Range("I14").Value = '28 minutes late'

$ws.Range("A15").Value = 'Anthony Martin'
$ws.Range("B15").Value = 'TRK-1014'
$ws.Range("C15").Value = 'West Plano Project'
$ws.Range("D15").Value = '07:00 AM'
$ws.Range("E15").Value = '04:00 PM'
$ws.Range("F15").Value = '07:16 AM'
$ws.Range("G15").Value = '03:58 PM'
$ws.Range("H15").Value = 'Late'
$ws.Range("I15").Value = '16 minutes late'

$ws.Range("A16").Value = 'Charles Anderson'
$ws.Range("B16").Value = 'TRK-1010'
$ws.Range("C16").Value = 'West Plano Project'
$ws.Range("D16").Value = '06:00 AM'
$ws.Range("E16").Value = '03:00 PM'
$ws.Range("F16").Value = '06:05 AM'
$ws.Range("G16").Value = '02:48 PM'
$ws.Range("H16").Value = 'On Time'
$ws.Range("I16").Value = ''
